$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Data" to "Summary"
$ws.Name = "Summary"

# Shift the existing table (rows 5-9) down by 4 rows, making room for a new
# "Source Type" line at row 7. After the insert:
#   old row 5 -> 9   (Micro / SMEs / MSMEs header)
#   old row 6 -> 10  (Enterprises (absolute #))
#   old row 7 -> 11  (Enterprises density (per 1000 people))
#   old row 8 -> 12  (Enterprises (% of total))
#   old row 9 -> 13  (Source: ISTEEBU, 2010)
$ws.Rows("5:8").Insert()

# New bold+underline named style used for the "Source Type" line
$titleStyle = $wb.Styles.Add("title_")
$titleStyle.Font.Bold = $true
$titleStyle.Font.Underline = $true

# New "Source Type" row (bold + underline)
$ws.Range("A7").Value = "Source Type: Statistical Institution (Most Widely Used)"
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Underline = $true

# New citation rows at the bottom of the sheet
$ws.Range("A21").Value = "ISTEEBU"
$ws.Range("A21").Font.Bold = $true

$ws.Range("A22").Value = "L'Institut de Statistiques et d’Etudes Economiques du Burundi (ISTEEBU), ""ANNUAIRE STATISTIQUE DU BURUNDI 2011"", Répartition des entreprises selon le nombre de travailleurs actifs (au 31 décembre), 2013, p. 160. Available at http://www.isteebu.bi/images/annuaires/annuaire%202011.pdf%20vf.pdf"
$ws.Range("A22").Font.Italic = $true
